$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Collapse the three detailed "CORE COMPETENCIES" paragraphs into a
# single summary paragraph listing just the three category names.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count

$firstIdx = -1
$lastIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $txt = $paras.Item($i).Range.Text
    if ($txt -like "Data Visualization & Design:*Interactive Dashboards*") {
        $firstIdx = $i
    }
    if ($txt -like "Technical Visualization:*Programming:*") {
        $lastIdx = $i
    }
}

if ($firstIdx -gt 0 -and $lastIdx -ge $firstIdx) {
    $startRange = $paras.Item($firstIdx).Range
    $endRange = $paras.Item($lastIdx).Range
    $competencyRange = $d.Range($startRange.Start, $endRange.End)
    $bullet = [char]0x2022
    $competencyRange.Text = "Data Visualization & Design " + $bullet + " Geospatial Analysis & Mapping " + $bullet + " Technical Visualization`r"
}

# ---------------------------------------------------------------------------
# Change 2: Add a new "TECHNICAL SKILLS" section (heading + three summary
# lines) right before the closing "For a more detailed..." paragraph.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count

$ledIdx = -1
for ($i = 1; $i -le $count; $i++) {
    if ($paras.Item($i).Range.Text -like "*Led multi-million dollar research projects*") {
        $ledIdx = $i
    }
}

if ($ledIdx -gt 0) {
    $anchor = $paras.Item($ledIdx).Range
    $anchor.InsertParagraphAfter()

    $skillLines = @(
        "TECHNICAL SKILLS",
        "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design",
        "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing",
        "TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing"
    )

    # First new paragraph becomes the "TECHNICAL SKILLS" Heading2 line.
    $newHeadingPara = $d.Paragraphs.Item($ledIdx + 1)
    $newHeadingPara.Range.Text = $skillLines[0]
    $newHeadingPara.Style = $d.Styles.Item("Heading 2")

    # Remaining lines are plain (Normal-style) paragraphs following the heading.
    $prevParaIdx = $ledIdx + 1
    for ($j = 1; $j -lt $skillLines.Length; $j++) {
        $prevPara = $d.Paragraphs.Item($prevParaIdx)
        $prevPara.Range.InsertParagraphAfter()
        $thisParaIdx = $prevParaIdx + 1
        $thisPara = $d.Paragraphs.Item($thisParaIdx)
        $thisPara.Style = $d.Styles.Item("Normal")
        $thisPara.Range.Text = $skillLines[$j]
        $prevParaIdx = $thisParaIdx
    }
}

Write-Host "Done"
